$d = $word.ActiveDocument

function Add-FigRef {
    param(
        [string]$OldText,
        [string[]]$Refs
    )

    $rng = $d.Content
    $found = $rng.Find.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND: $OldText"
        return
    }

    # Remove the trailing colon ':' from the matched sentence, then
    # rebuild "(рис. [N?], [M?]...):" in its place.
    $trimmedEnd = $rng.End - 1
    $colonRange = $d.Range($trimmedEnd, $rng.End)
    $colonRange.Text = ""

    $insPoint = $trimmedEnd
    $tailRange = $d.Range($insPoint, $insPoint)
    $tailRange.InsertAfter("(рис.")
    $insPoint = $insPoint + 5

    $first = $true
    foreach ($ref in $Refs) {
        $sep = $d.Range($insPoint, $insPoint)
        if ($first) {
            $sep.InsertAfter(" ")
            $insPoint = $insPoint + 1
            $first = $false
        } else {
            $sep.InsertAfter(", ")
            $insPoint = $insPoint + 2
        }

        $openBr = $d.Range($insPoint, $insPoint)
        $openBr.InsertAfter("[")
        $insPoint = $insPoint + 1

        $numStart = $insPoint
        $numIns = $d.Range($insPoint, $insPoint)
        $numIns.InsertAfter($ref)
        $insPoint = $insPoint + $ref.Length
        $numRange = $d.Range($numStart, $insPoint)
        $numRange.Bold = 1

        $closeBr = $d.Range($insPoint, $insPoint)
        $closeBr.InsertAfter("]")
        $insPoint = $insPoint + 1
    }

    $tail2 = $d.Range($insPoint, $insPoint)
    $tail2.InsertAfter("):")
    $insPoint = $insPoint + 2
}

Add-FigRef "Настроим git для первоначальной работы:" @("1?")
Add-FigRef "Используем команду ssh-keygen, чтобы сгененерировать публичный и приватный ssh-ключи, которые понадобятся для аутентификации:" @("2?")
Add-FigRef "Скопировав содержимое файла с ключем, добавим его на Github:" @("3?")
Add-FigRef "Используем указанный репозиторий в качестве шаблона и убедимся, что клонирование прошло успешно:" @("4?")
Add-FigRef "Создадим рабочий каталог и подкаталоги в домашнем каталоге, в которые мы впоследствии будем клонировать наш репозиторий, и перейдём в него:" @("5?")
Add-FigRef "Убедимся, что клонирование репозитория прошло успешно:" @("6?")
Add-FigRef "Введём последовательность команд для создания правильной файловой структуры, и убедимся в её наличии:" @("7?")
Add-FigRef "Добавим изменения в репозиторий командами git, после чего отправим репозиторий на Github:" @("8?", "9?")
Add-FigRef "Убедимся, что изменения кода прошли на Github:" @("10?")
